# FAKE_DATA_2019.xlsx - "created sip with 4 digits to export in PANON and Address"
#
# - remove the EGID / EWID columns (N:O) entirely, shifting VERMÖGEN,
#   STEUERBARESEINKOMMEN, HASEL, HASSH, AMOUNT two columns to the left
# - remove the three extra duplicate/test records in rows 9-11
# - fix ZIVILSTAND for row 3 (Brunner, Berta): "Civil stat#L" -> "Civil stat#H"
# - fix HAUSNR for row 4 (Christen, Claudio): "3c!" -> "3c" (4 digits -> plain)
# - narrow column B now that the longer EGID/EWID-era content is gone
# - move the active selection to R9 (last used cell after the cleanup)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the EGID (N) and EWID (O) columns.
$ws.Columns("N:O").Delete()

# Drop the synthetic extra rows (Anton / Anita / Ali).
$ws.Rows("9:11").Delete()

# Data corrections.
$ws.Range("F3").Value = "Civil stat#H"
$ws.Range("K4").Value = "3c"

# Column B no longer needs to be as wide.
$ws.Columns("B").ColumnWidth = 10.83

# Update the selected cell to match the new extent of the data.
[void]$ws.Range("R9").Select()
